$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the stimulus images (column B) for the new "recall_cross" set.
$ws.Range("B2").Value = "common/recall_cross_EL.png"
$ws.Range("B3").Value = "common/recall_cross_ER.png"
$ws.Range("B4").Value = "common/recall_cross_IL.png"
$ws.Range("B5").Value = "common/recall_cross_IR.png"
$ws.Range("B6").Value = "common/recall_cross_IT.png"

# The "ori" (orientation) column is no longer used - clear it out.
$ws.Range("C1:C6").ClearContents()

# C2 keeps a (now empty) formatted cell behind - touch its format so a
# style record for it survives the save.
$ws.Range("C2").Locked = $true

# Widen column B so the longer image paths are readable.
$ws.Columns.Item(2).ColumnWidth = 54.42

# Page setup tweaks recorded alongside the data edits.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the selection where the author left it.
$ws.Range("C2").Select() | Out-Null
